$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Price (column D) updates - forced to remain text
Set-TextValue $ws 'D2' '68.807.42'
Set-TextValue $ws 'D3' '2.436.42'
Set-TextValue $ws 'D5' '559.65'
Set-TextValue $ws 'D6' '162.03'
Set-TextValue $ws 'D11' '0.331'
Set-TextValue $ws 'D14' '68.702.48'
Set-TextValue $ws 'D15' '2.886.20'
Set-TextValue $ws 'D16' '23.26'
Set-TextValue $ws 'D17' '2.438.44'
Set-TextValue $ws 'D18' '10.54'
Set-TextValue $ws 'D19' '339.06'
Set-TextValue $ws 'D20' '6.93'
Set-TextValue $ws 'D22' '1.92'
Set-TextValue $ws 'D24' '67.01'
Set-TextValue $ws 'D25' '3.71'
Set-TextValue $ws 'D26' '2.566.55'
Set-TextValue $ws 'D27' '1.01'
Set-TextValue $ws 'D28' '8.20'
Set-TextValue $ws 'D29' '0.0₃0820'
Set-TextValue $ws 'D30' '7.14'
Set-TextValue $ws 'D31' '0.999'
Set-TextValue $ws 'D33' '427.90'
Set-TextValue $ws 'D35' '159.12'
Set-TextValue $ws 'D36' '18.99'
Set-TextValue $ws 'D38' '17.99'
Set-TextValue $ws 'D42' '4.35'
Set-TextValue $ws 'D45' '3.34'
Set-TextValue $ws 'D46' '130.33'
Set-TextValue $ws 'D47' '0.0718'
Set-TextValue $ws 'D48' '0.480'

# Volume(1h) (column E) updates
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('E9').Value = '  +10.78%  '
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('E12').Value = '  -5.89%  '
$ws.Range('E13').Value = '  +5.30%  '
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('E22').Value = '  +2.45%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('E25').Value = '  +2.43%  '
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +2.88%  '
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('E41').Value = '  +3.72%  '
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('E51').Value = '  +1.16%  '
